$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, matching the source
# workbook (values like "0.999" or "63.253.90" must not be coerced to numbers).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.253.90'
$ws.Range('E2').Value = '  -3.70%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.596.25'
$ws.Range('E3').Value = '  -2.44%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.28'
$ws.Range('E5').Value = '  -4.93%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.55'
$ws.Range('E6').Value = '  -3.84%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('E8').Value = '  -4.70%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.593.92'
$ws.Range('E9').Value = '  -2.35%  '

$ws.Range('E10').Value = '  -8.40%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.72'
$ws.Range('E11').Value = '  -2.19%  '

$ws.Range('E12').Value = '  -0.33%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.373'
$ws.Range('E13').Value = '  -6.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.68'
$ws.Range('E14').Value = '  -4.84%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.055.21'
$ws.Range('E15').Value = '  -2.50%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  -9.21%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.105.44'
$ws.Range('E17').Value = '  -3.75%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.612.09'
$ws.Range('E18').Value = '  -1.46%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.83'
$ws.Range('E19').Value = '  -5.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.39'
$ws.Range('E20').Value = '  -0.86%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.43'
$ws.Range('E21').Value = '  -7.04%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '337.43'
$ws.Range('E22').Value = '  -4.46%  '

$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.68'
$ws.Range('E24').Value = '  -3.80%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.78'
$ws.Range('E25').Value = '  +0.87%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000105'
$ws.Range('E26').Value = '  -7.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '576.33'
$ws.Range('E27').Value = '  +2.46%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.99'
$ws.Range('E28').Value = '  -6.09%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.53'
$ws.Range('E29').Value = '  -5.72%  '

$ws.Range('E30').Value = '  +0.48%  '

$ws.Range('E31').Value = '  -2.92%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.67'
$ws.Range('E32').Value = '  -5.19%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.02'
$ws.Range('E33').Value = '  -5.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.69'
$ws.Range('E34').Value = '  -6.35%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.44'
$ws.Range('E35').Value = '  -3.58%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.28'
$ws.Range('E36').Value = '  -3.88%  '

$ws.Range('E37').Value = '  -0.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.396'
$ws.Range('E38').Value = '  -6.01%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.46'
$ws.Range('E39').Value = '  -4.91%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '154.01'
$ws.Range('E40').Value = '  +0.67%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.84'
$ws.Range('E41').Value = '  -6.33%  '

$ws.Range('E42').Value = '  -0.08%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.33'
$ws.Range('E43').Value = '  -3.59%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.47'
$ws.Range('E44').Value = '  +0.64%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '156.56'
$ws.Range('E45').Value = '  -3.12%  '

$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.81'
$ws.Range('E46').Value = '  -6.84%  '

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.88'
$ws.Range('E47').Value = '  -1.71%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0575'
$ws.Range('E48').Value = '  -6.37%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.623'
$ws.Range('E49').Value = '  -3.16%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0985'
$ws.Range('E50').Value = '  -3.11%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0243'
$ws.Range('E51').Value = '  -5.66%  '
